# "Generate Report for Handback"
#
# The localization hand-off for the two e2e markdown files has come back
# from the translators, so this report now reflects a completed handback:
#   * Status flips from "Ready for handoff" to "Handed back: in sync with en-US"
#     (shown on the Overview sheet as well as each language sheet's Status column).
#   * Each language sheet gets its "Latest Target File" / "Latest Handback File"
#     columns populated (with a hyperlink on the target-file cell, matching the
#     existing source-file hyperlink style) and a fresh "Latest Handback DateTime".

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0c727d2f8499be35219865b79759b978f7ac3c4a/e2e/"

$mdFile1 = "39b05a38-cb78-449e-9597-5c512e044c2c.md"
$mdFile2 = "b73dfcf0-6ec0-4e0b-9b79-830ad2e0e8ce.md"

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: the zh-cn / de-de status columns (E & F) move from
# "Ready for handoff" to "Handed back: in sync with en-US" for both rows.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Columns.Item(5).EntireColumn.AutoFit()
$wsOverview.Columns.Item(6).EntireColumn.AutoFit()

# ---------------------------------------------------------------------------
# zh-cn sheet: Status column + newly-populated target/handback columns.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Range("I2").Value = $mdFile1
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), ($baseUrl + $mdFile1), "", "", $mdFile1) | Out-Null
$wsZh.Range("J2").Value = "39b05a38-cb78-449e-9597-5c512e044c2c.4e4aa8e48aa9029afbc2e7c492c2977b72fb3b8d.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-06 07:12:53"

$wsZh.Range("I3").Value = $mdFile2
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), ($baseUrl + $mdFile2), "", "", $mdFile2) | Out-Null
$wsZh.Range("J3").Value = "b73dfcf0-6ec0-4e0b-9b79-830ad2e0e8ce.27c19c26c70fb024e3fd9a90fcffe9787a6ada25.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-09-06 07:12:53"

$wsZh.Columns.Item(3).EntireColumn.AutoFit()
$wsZh.Columns.Item(9).EntireColumn.AutoFit()
$wsZh.Columns.Item(10).EntireColumn.AutoFit()

# ---------------------------------------------------------------------------
# de-de sheet: Status column + newly-populated target/handback columns.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Range("I2").Value = $mdFile1
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), ($baseUrl + $mdFile1), "", "", $mdFile1) | Out-Null
$wsDe.Range("J2").Value = "39b05a38-cb78-449e-9597-5c512e044c2c.4e4aa8e48aa9029afbc2e7c492c2977b72fb3b8d.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-06 07:13:01"

$wsDe.Range("I3").Value = $mdFile2
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), ($baseUrl + $mdFile2), "", "", $mdFile2) | Out-Null
$wsDe.Range("J3").Value = "b73dfcf0-6ec0-4e0b-9b79-830ad2e0e8ce.27c19c26c70fb024e3fd9a90fcffe9787a6ada25.de-de.xlf"
$wsDe.Range("K3").Value = "2016-09-06 07:13:01"

$wsDe.Columns.Item(3).EntireColumn.AutoFit()
$wsDe.Columns.Item(9).EntireColumn.AutoFit()
$wsDe.Columns.Item(10).EntireColumn.AutoFit()

Write-Host "Handback report generated."
